$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 11.20127533333333
$ws.Range("H2").Value = 33.603826
$ws.Range("I2").Value = 0.1186573945858706
$ws.Range("J2").Value = 0.1186573945858706
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 153.5290173333333
$ws.Range("N2").Value = 460.587052
$ws.Range("O2").Value = 0.3172206968818489
$ws.Range("P2").Value = 0.317220696881849
$ws.Range("Q2").Value = 1719.720794806772
$ws.Range("R2").Value = 15477.48715326095
$ws.Range("S2").Value = 0.03764058140071441
$ws.Range("T2").Value = 0.03764058140071441

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 11.20127533333333
$ws.Range("H3").Value = 33.603826
$ws.Range("I3").Value = 0.1186573945858706
$ws.Range("J3").Value = 0.1186573945858706
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 168.7997026666667
$ws.Range("N3").Value = 506.3991080000001
$ws.Range("O3").Value = 0.3487728915577651
$ws.Range("P3").Value = 0.3487728915577651
$ws.Range("Q3").Value = 1890.771945754135
$ws.Range("R3").Value = 17016.94751178721
$ws.Range("S3").Value = 0.0413844826144248
$ws.Range("T3").Value = 0.0413844826144248

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 11.20127533333333
$ws.Range("H4").Value = 33.603826
$ws.Range("I4").Value = 0.1186573945858706
$ws.Range("J4").Value = 0.1186573945858706
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 68.09032333333333
$ws.Range("N4").Value = 204.27097
$ws.Range("O4").Value = 0.1406878008722904
$ws.Range("P4").Value = 0.1406878008722904
$ws.Range("Q4").Value = 762.6984591923577
$ws.Range("R4").Value = 6864.286132731219
$ws.Range("S4").Value = 0.01669364790152176
$ws.Range("T4").Value = 0.01669364790152176

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 11.20127533333333
$ws.Range("H5").Value = 33.603826
$ws.Range("I5").Value = 0.1186573945858706
$ws.Range("J5").Value = 0.1186573945858706
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 93.562673
$ws.Range("N5").Value = 280.688019
$ws.Range("O5").Value = 0.1933186106880956
$ws.Range("P5").Value = 0.1933186106880956
$ws.Range("Q5").Value = 1048.021261195633
$ws.Range("R5").Value = 9432.191350760693
$ws.Range("S5").Value = 0.02293868266920967
$ws.Range("T5").Value = 0.02293868266920967

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 47.94465366666667
$ws.Range("H6").Value = 143.833961
$ws.Range("I6").Value = 0.5078874966566524
$ws.Range("J6").Value = 0.5078874966566524
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 153.5290173333333
$ws.Range("N6").Value = 460.587052
$ws.Range("O6").Value = 0.3172206968818489
$ws.Range("P6").Value = 0.317220696881849
$ws.Range("Q6").Value = 7360.895563830331
$ws.Range("R6").Value = 66248.06007447297
$ws.Range("S6").Value = 0.161112425627001
$ws.Range("T6").Value = 0.161112425627001

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 47.94465366666667
$ws.Range("H7").Value = 143.833961
$ws.Range("I7").Value = 0.5078874966566524
$ws.Range("J7").Value = 0.5078874966566524
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 168.7997026666667
$ws.Range("N7").Value = 506.3991080000001
$ws.Range("O7").Value = 0.3487728915577651
$ws.Range("P7").Value = 0.3487728915577651
$ws.Range("Q7").Value = 8093.043283389646
$ws.Range("R7").Value = 72837.38955050681
$ws.Range("S7").Value = 0.1771373907949754
$ws.Range("T7").Value = 0.1771373907949754

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 47.94465366666667
$ws.Range("H8").Value = 143.833961
$ws.Range("I8").Value = 0.5078874966566524
$ws.Range("J8").Value = 0.5078874966566524
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 68.09032333333333
$ws.Range("N8").Value = 204.27097
$ws.Range("O8").Value = 0.1406878008722904
$ws.Range("P8").Value = 0.1406878008722904
$ws.Range("Q8").Value = 3264.566970268019
$ws.Range("R8").Value = 29381.10273241218
$ws.Range("S8").Value = 0.07145357499515717
$ws.Range("T8").Value = 0.07145357499515718

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 47.94465366666667
$ws.Range("H9").Value = 143.833961
$ws.Range("I9").Value = 0.5078874966566524
$ws.Range("J9").Value = 0.5078874966566524
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 93.562673
$ws.Range("N9").Value = 280.688019
$ws.Range("O9").Value = 0.1933186106880956
$ws.Range("P9").Value = 0.1933186106880956
$ws.Range("Q9").Value = 4485.829953112585
$ws.Range("R9").Value = 40372.46957801326
$ws.Range("S9").Value = 0.09818410523951887
$ws.Range("T9").Value = 0.09818410523951887

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 25.59984766666667
$ws.Range("H10").Value = 76.799543
$ws.Range("I10").Value = 0.2711844085184091
$ws.Range("J10").Value = 0.2711844085184091
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 153.5290173333333
$ws.Range("N10").Value = 460.587052
$ws.Range("O10").Value = 0.3172206968818489
$ws.Range("P10").Value = 0.317220696881849
$ws.Range("Q10").Value = 3930.319456146359
$ws.Range("R10").Value = 35372.87510531723
$ws.Range("S10").Value = 0.08602530705370175
$ws.Range("T10").Value = 0.08602530705370176

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 25.59984766666667
$ws.Range("H11").Value = 76.799543
$ws.Range("I11").Value = 0.2711844085184091
$ws.Range("J11").Value = 0.2711844085184091
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 168.7997026666667
$ws.Range("N11").Value = 506.3991080000001
$ws.Range("O11").Value = 0.3487728915577651
$ws.Range("P11").Value = 0.3487728915577651
$ws.Range("Q11").Value = 4321.246674445294
$ws.Range("R11").Value = 38891.22007000765
$ws.Range("S11").Value = 0.09458177030434778
$ws.Range("T11").Value = 0.09458177030434778

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 25.59984766666667
$ws.Range("H12").Value = 76.799543
$ws.Range("I12").Value = 0.2711844085184091
$ws.Range("J12").Value = 0.2711844085184091
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 68.09032333333333
$ws.Range("N12").Value = 204.27097
$ws.Range("O12").Value = 0.1406878008722904
$ws.Range("P12").Value = 0.1406878008722904
$ws.Range("Q12").Value = 1743.101904907412
$ws.Range("R12").Value = 15687.91714416671
$ws.Range("S12").Value = 0.03815233806530779
$ws.Range("T12").Value = 0.0381523380653078

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 25.59984766666667
$ws.Range("H13").Value = 76.799543
$ws.Range("I13").Value = 0.2711844085184091
$ws.Range("J13").Value = 0.2711844085184091
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 93.562673
$ws.Range("N13").Value = 280.688019
$ws.Range("O13").Value = 0.1933186106880956
$ws.Range("P13").Value = 0.1933186106880956
$ws.Range("Q13").Value = 2395.190176086146
$ws.Range("R13").Value = 21556.71158477532
$ws.Range("S13").Value = 0.05242499309505182
$ws.Range("T13").Value = 0.05242499309505182

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 9.654369000000001
$ws.Range("H14").Value = 28.963107
$ws.Range("I14").Value = 0.1022707002390678
$ws.Range("J14").Value = 0.1022707002390678
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 153.5290173333333
$ws.Range("N14").Value = 460.587052
$ws.Range("O14").Value = 0.3172206968818489
$ws.Range("P14").Value = 0.317220696881849
$ws.Range("Q14").Value = 1482.225785543396
$ws.Range("R14").Value = 13340.03206989056
$ws.Range("S14").Value = 0.03244238280043175
$ws.Range("T14").Value = 0.03244238280043175

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 9.654369000000001
$ws.Range("H15").Value = 28.963107
$ws.Range("I15").Value = 0.1022707002390678
$ws.Range("J15").Value = 0.1022707002390678
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 168.7997026666667
$ws.Range("N15").Value = 506.3991080000001
$ws.Range("O15").Value = 0.3487728915577651
$ws.Range("P15").Value = 0.3487728915577651
$ws.Range("Q15").Value = 1629.654616634284
$ws.Range("R15").Value = 14666.89154970856
$ws.Range("S15").Value = 0.03566924784401709
$ws.Range("T15").Value = 0.03566924784401709

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 9.654369000000001
$ws.Range("H16").Value = 28.963107
$ws.Range("I16").Value = 0.1022707002390678
$ws.Range("J16").Value = 0.1022707002390678
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 68.09032333333333
$ws.Range("N16").Value = 204.27097
$ws.Range("O16").Value = 0.1406878008722904
$ws.Range("P16").Value = 0.1406878008722904
$ws.Range("Q16").Value = 657.3691067893101
$ws.Range("R16").Value = 5916.32196110379
$ws.Range("S16").Value = 0.01438823991030367
$ws.Range("T16").Value = 0.01438823991030367

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 9.654369000000001
$ws.Range("H17").Value = 28.963107
$ws.Range("I17").Value = 0.1022707002390678
$ws.Range("J17").Value = 0.1022707002390678
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 93.562673
$ws.Range("N17").Value = 280.688019
$ws.Range("O17").Value = 0.1933186106880956
$ws.Range("P17").Value = 0.1933186106880956
$ws.Range("Q17").Value = 903.2885697683371
$ws.Range("R17").Value = 8129.597127915034
$ws.Range("S17").Value = 0.01977082968431527
$ws.Range("T17").Value = 0.01977082968431527
